# Sara-Alert-Format-Exposure-Workflow.xlsx
# "allow vaccine table to be populated on import"
#
# Adds a second vaccination block (Vaccine 1 / Vaccine 2, each with Group
# Name, Product Name, Administration Date, Dose Number, Notes) to the
# Monitorees sheet starting at column CY (103), and populates sample data
# for the first few rows, the way SaraAlert's import template fixture does.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Column headers (row 1)
# ---------------------------------------------------------------------
$headers = @{
    103 = "Vaccine 1 Group Name";
    104 = "Vaccine 1 Product Name";
    105 = "Vaccine 1 Administration Date";
    106 = "Vaccine 1 Dose Number";
    107 = "Vaccine 1 Notes";
    108 = "Vaccine 2 Group Name";
    109 = "Vaccine 2 Product Name";
    110 = "Vaccine 2 Administration Date";
    111 = "Vaccine 2 Dose Number";
    112 = "Vaccine 2 Notes";
}

# The two "Administration Date" columns are stored as text (they keep the
# literal yyyy-mm-dd string instead of being converted to a date serial),
# matching the rest of the workbook's "text" formatted columns.
$textDateCols = @(105, 110)

foreach ($col in $textDateCols) {
    $ws.Range($ws.Cells.Item(1, $col), $ws.Cells.Item(12, $col)).NumberFormat = "@"
}

foreach ($col in $headers.Keys) {
    $ws.Cells.Item(1, $col).Value = $headers[$col]
}

# ---------------------------------------------------------------------
# Sample data (rows 2-7)
# ---------------------------------------------------------------------
$rowData = @{
    2  = @{ 103 = "COVID-19"; 104 = "Moderna COVID-19 Vaccine";         105 = "2020-06-01"; 106 = 1; 107 = "notes 1";
            108 = "COVID-19"; 109 = "Moderna COVID-19 Vaccine";         110 = "2020-06-20"; 111 = 2; 112 = "notes 2" };
    3  = @{ 103 = "COVID-19"; 104 = "Pfizer-BioNTech COVID-19 Vaccine"; 105 = "2020-06-02"; 106 = 1;
            108 = "COVID-19"; 109 = "Pfizer-BioNTech COVID-19 Vaccine"; 110 = "2020-06-21"; 111 = 2 };
    4  = @{ 103 = "COVID-19"; 104 = "Unknown";                          105 = "2020-06-04"; 106 = 1;
            108 = "COVID-19"; 109 = "Unknown";                          110 = "2020-06-22"; 111 = 2 };
    5  = @{ 103 = "COVID-19"; 104 = "Moderna COVID-19 Vaccine";         105 = "2020-06-01"; 106 = 1 };
    6  = @{ 103 = "COVID-19"; 104 = "Janssen (J&J) COVID-19 Vaccine";   105 = "2020-06-03"; 106 = 1 };
    7  = @{ 103 = "COVID-19"; 104 = "Unknown";                          105 = "2020-06-02"; 106 = 1 };
}

foreach ($r in $rowData.Keys) {
    $cells = $rowData[$r]
    foreach ($col in $cells.Keys) {
        $ws.Cells.Item($r, $col).Value = $cells[$col]
    }
}

# ---------------------------------------------------------------------
# Column widths for the new columns (best match for the bestFit widths
# Excel would have computed for this content)
# ---------------------------------------------------------------------
$colWidths = @{
    103 = 20.33203125;
    104 = 31;
    105 = 25.6640625;
    106 = 21.1640625;
    107 = 14.5;
    108 = 20.33203125;
    109 = 31;
    110 = 25.6640625;
    111 = 21.1640625;
    112 = 14.5;
}

foreach ($col in $colWidths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $colWidths[$col] - (5.0 / 7.0)
}

# ---------------------------------------------------------------------
# Reset the view back to the top-left corner / A1, as is typical after
# the sheet has been re-saved from a fresh load (the previous revision
# had scrolled over to show the far-right columns).
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
